$d = $word.ActiveDocument

# The paragraph currently reads (across three differently-formatted runs):
#   "Peta di dokumen ini disiapkan oleh Jenik Hollan, CzechGlobe ("
#   "http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2018/"   (hyperlink-styled)
#   ")."
# It needs to become a single plain run with the year bumped to 2022:
#   "Peta di dokumen ini disiapkan oleh Jenik Hollan, CzechGlobe (http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/)."

$f = $d.Content.Find
$found = $f.Execute("Peta di dokumen ini disiapkan oleh Jenik Hollan, CzechGlobe (http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2018/).", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $r = $f.Parent
    $r.Delete()
    $r.InsertAfter("Peta di dokumen ini disiapkan oleh Jenik Hollan, CzechGlobe (http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/).")
}
